# The "Rules" sheet has a small lookup table in B3:E11. The last data
# row (row 11) originally carried the label "R40" in column B (as text,
# since it is a shared string alongside the other row labels R10/R20/R30).
# This edit changes that label to the text "1".
#
# A leading apostrophe is used so Excel stores the value as literal text
# (shared string) rather than re-interpreting the numeric-looking "1" as
# a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B11").Value = "'1"
